$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, matching the style of the existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# New data value in H2
$ws.Range("H2").Value = 1
